# Update to calculate scores relative to median AFTER merging with zip/census tract.
#
# Column C ("Median Value") previously held each school's poverty ratio
# divided by the pre-merge median; it is now recomputed against the median
# produced after merging in the zip/census-tract data (every ratio below is
# the old ratio divided by the old median, 1.185821697099892 - which is why
# row 7, the median row itself, now reads exactly 1). Column D ("Tier") is
# refreshed to match the re-ranked quartiles of the new Column C values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.028985507246377
$ws.Range("D2").Value = "4th Tier"
$ws.Range("C3").Value = 0.5090579710144928
$ws.Range("C4").Value = 1.957971014492754
$ws.Range("C5").Value = 1.318840579710145
$ws.Range("D5").Value = "3rd Tier"
$ws.Range("C6").Value = 1.540760869565217
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 0.6898550724637681
$ws.Range("C9").Value = 1.855676328502415
$ws.Range("C10").Value = 0.7608695652173912
$ws.Range("C11").Value = 1.00054347826087
$ws.Range("C12").Value = 1.032608695652174
$ws.Range("D12").Value = "4th Tier"
$ws.Range("C13").Value = 1.087409420289855
$ws.Range("D13").Value = "4th Tier"
$ws.Range("C14").Value = 1.334692028985507
$ws.Range("C15").Value = 1.458333333333333
$ws.Range("D15").Value = "2nd Tier"
$ws.Range("C16").Value = 0.8293478260869566
$ws.Range("C17").Value = 0.9146286231884058
$ws.Range("D17").Value = "Below Median"
$ws.Range("C18").Value = 0.5727657004830917
$ws.Range("C19").Value = 0.8510466988727858
$ws.Range("D19").Value = "Below Median"
$ws.Range("C20").Value = 0.7059178743961352
$ws.Range("C21").Value = 0.6518115942028985
$ws.Range("C22").Value = 0.5217391304347826
$ws.Range("C23").Value = 0.6105072463768116
$ws.Range("C24").Value = 1.389855072463768
$ws.Range("C25").Value = 1.675724637681159
$ws.Range("C26").Value = 1.361111111111111
$ws.Range("C27").Value = 1.499547101449275
$ws.Range("C28").Value = 1.001811594202898
$ws.Range("C29").Value = 0.5757246376811593
$ws.Range("C30").Value = 1.43677536231884
$ws.Range("C31").Value = 1.059581320450886
$ws.Range("D31").Value = "4th Tier"
$ws.Range("C32").Value = 1.678985507246377
$ws.Range("C33").Value = 1.016606280193237
$ws.Range("D33").Value = "4th Tier"
$ws.Range("C34").Value = 0.9710144927536231
$ws.Range("D34").Value = "Below Median"
$ws.Range("C35").Value = 0.4839975845410627
$ws.Range("C36").Value = 1.356884057971014
$ws.Range("C37").Value = 0.7355072463768115
$ws.Range("C38").Value = 1.27536231884058
$ws.Range("D38").Value = "3rd Tier"
$ws.Range("C39").Value = 1.306159420289855
$ws.Range("D39").Value = "3rd Tier"
$ws.Range("C40").Value = 1.151449275362319
$ws.Range("C41").Value = 1.109601449275362
$ws.Range("D41").Value = "4th Tier"
$ws.Range("C42").Value = 0.8327294685990339
$ws.Range("C43").Value = 1.534420289855072
$ws.Range("C44").Value = 0.8876811594202898
$ws.Range("D44").Value = "Below Median"
$ws.Range("C45").Value = 1.113405797101449
$ws.Range("C46").Value = 0.601086956521739
$ws.Range("C47").Value = 0.9658816425120773
$ws.Range("D47").Value = "Below Median"
$ws.Range("C48").Value = 1.185688405797101
$ws.Range("C49").Value = 1.41268115942029
$ws.Range("C50").Value = 1.071557971014493
$ws.Range("D50").Value = "4th Tier"
$ws.Range("C51").Value = 0.8834541062801933
$ws.Range("D51").Value = "Below Median"
$ws.Range("C52").Value = 0.6644927536231884
$ws.Range("C53").Value = 1.204710144927536
$ws.Range("C54").Value = 0.9035326086956521
$ws.Range("D54").Value = "Below Median"
$ws.Range("C55").Value = 1.123188405797101
$ws.Range("C56").Value = 0.9184782608695652
$ws.Range("D56").Value = "Below Median"
$ws.Range("C57").Value = 0.527536231884058
$ws.Range("C58").Value = 0.6124999999999999
$ws.Range("C59").Value = 0.3043478260869565
$ws.Range("C60").Value = 0.5126811594202898
$ws.Range("C61").Value = 0.8211050724637681
$ws.Range("C62").Value = 1.341032608695652
$ws.Range("C63").Value = 0.6391304347826087
$ws.Range("C64").Value = 0.5807971014492753
$ws.Range("C65").Value = 0.6625905797101449
$ws.Range("C66").Value = 0.3719806763285024
$ws.Range("C67").Value = 0.9739130434782608
$ws.Range("D67").Value = "Below Median"
$ws.Range("C68").Value = 1.911684782608696
$ws.Range("C69").Value = 0.6036231884057971
$ws.Range("C70").Value = 0.8876811594202898
$ws.Range("D70").Value = "Below Median"
$ws.Range("C71").Value = 1.607971014492753
$ws.Range("C72").Value = 1.557246376811594
$ws.Range("C73").Value = 0.5427536231884057
$ws.Range("C74").Value = 1.327898550724637
$ws.Range("D74").Value = "3rd Tier"
$ws.Range("C75").Value = 1.296014492753623
$ws.Range("D75").Value = "3rd Tier"
$ws.Range("C76").Value = 1.348429951690821
$ws.Range("C77").Value = 0.7684782608695652
$ws.Range("C78").Value = 1.480525362318841
$ws.Range("D78").Value = "2nd Tier"
$ws.Range("C79").Value = 0.6657608695652174
$ws.Range("C80").Value = 0.6467391304347826
